$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows: Indice, Distancia, max, min, Tempo
$data = @(
    @(0, 7412.066666666667, 8011, 6888, 0.1148659467697144),
    @(1, 7914.2,            8561, 7472, 0.1208336750666301),
    @(2, 8134.533333333334, 8709, 7353, 0.1223240375518799),
    @(3, 8661.866666666667, 9485, 7842, 0.1191922346750895),
    @(4, 7755.366666666667, 8157, 7334, 0.1149209499359131),
    @(5, 7961.333333333333, 8621, 7428, 0.1212686618169149),
    @(6, 7914.7,            8429, 7129, 0.1218839248021444),
    @(7, 7682.966666666666, 8354, 6682, 0.1208723942438761),
    @(8, 7146.166666666667, 7599, 6387, 0.1165302356084188),
    @(9, 7523.2,            7979, 6389, 0.1162815014521281)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
